$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New chronological order of "Periodo Mora" labels for rows 16-41
$periodos = @("2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112","2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212","2301","2302")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]   # E: Periodo Mora
    $ws.Cells.Item($row, 7).Value = 2450000          # G: Salario Basico
}

# F column ("Valor Mora"): every row is 70000 except the last data row (41),
# which keeps the odd legacy value (previously on the first row, 16)
for ($row = 16; $row -le 40; $row++) {
    $ws.Cells.Item($row, 6).Value = 70000
}
$ws.Cells.Item(41, 6).Value = 44333
